# Backlog.xlsx update as of 4/29
# Fill in "Confirmations" (I) = "Done" and "Actual Time" (K) values
# for the backlog items that have now been completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Backlog")

# Row 2 - item 9
$ws.Range("I2").Value = "Done"
$ws.Range("K2").Value = 10

# Row 9 - item 6
$ws.Range("I9").Value = "Done"
$ws.Range("K9").Value = 1

# Row 11 - item 12
$ws.Range("I11").Value = "Done"
$ws.Range("K11").Value = 4

# Row 12 - item 10
$ws.Range("I12").Value = "Done"
$ws.Range("K12").Value = 4

# Row 13 - item 11
$ws.Range("I13").Value = "Done"
$ws.Range("K13").Value = 2

# Row 15 - item 14
$ws.Range("I15").Value = "Done"
$ws.Range("K15").Value = 4

# Move the active selection to K15, matching where work left off
$ws.Activate()
$ws.Range("K15").Select()
